# Auto-generated: update Sheets via scheduled runner
# Applies refreshed market-price figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# to the per-job Phantom Profits tables across all 8 worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 395
$ws.Range("I31").Value = 395
$ws.Range("K31").Value = 1185
$ws.Range("M31").Value = -955
$ws.Range("H33").Value = 255.15384
$ws.Range("I33").Value = 255.15384
$ws.Range("K33").Value = 255.15384
$ws.Range("M33").Value = -26.15384
$ws.Range("H40").Value = 1174.25
$ws.Range("I40").Value = 1187.9333
$ws.Range("K40").Value = 1187.9333
$ws.Range("M40").Value = -1012.9333
$ws.Range("H43").Value = 1844.1111
$ws.Range("I43").Value = 1924.75
$ws.Range("K43").Value = 1924.75
$ws.Range("M43").Value = -1855.75
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
$ws.Range("H93").Value = 600999.5
$ws.Range("J93").Value = 600999.5
$ws.Range("L93").Value = 600999.5
$ws.Range("N93").Value = -605991.5
$ws.Range("H115").Value = 380
$ws.Range("I115").Value = 380
$ws.Range("K115").Value = 1140
$ws.Range("M115").Value = 427
$ws.Range("H138").Value = 1012.6875
$ws.Range("I138").Value = 880.93335
$ws.Range("J138").Value = 2989
$ws.Range("K138").Value = 2642.80005
$ws.Range("L138").Value = 8967
$ws.Range("M138").Value = 2497.19995
$ws.Range("N138").Value = -19247
$ws.Range("H141").Value = 3467.8948
$ws.Range("I141").Value = 3612.2778
$ws.Range("K141").Value = 10836.8334
$ws.Range("M141").Value = -5656.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3546.0625
$ws.Range("I61").Value = 3619.0908
$ws.Range("J61").Value = 3385.4
$ws.Range("K61").Value = 3619.0908
$ws.Range("L61").Value = 3385.4
$ws.Range("M61").Value = -3407.0908
$ws.Range("N61").Value = -3809.4
$ws.Range("H74").Value = 1699.579
$ws.Range("I74").Value = 1707.25
$ws.Range("J74").Value = 1686.4286
$ws.Range("K74").Value = 1707.25
$ws.Range("L74").Value = 1686.4286
$ws.Range("M74").Value = -833.25
$ws.Range("N74").Value = -3434.4286
$ws.Range("H77").Value = 1699.579
$ws.Range("I77").Value = 1707.25
$ws.Range("J77").Value = 1686.4286
$ws.Range("K77").Value = 8536.25
$ws.Range("L77").Value = 8432.143
$ws.Range("M77").Value = -4168.25
$ws.Range("N77").Value = -17168.143
$ws.Range("H88").Value = 2003
$ws.Range("I88").Value = 1696.5
$ws.Range("K88").Value = 1696.5
$ws.Range("M88").Value = -1290.5
$ws.Range("H91").Value = 2003
$ws.Range("I91").Value = 1696.5
$ws.Range("K91").Value = 1696.5
$ws.Range("M91").Value = -292.5
$ws.Range("H110").Value = 6587.6665
$ws.Range("I110").Value = 7687.375
$ws.Range("J110").Value = 4388.25
$ws.Range("K110").Value = 7687.375
$ws.Range("L110").Value = 4388.25
$ws.Range("M110").Value = -5642.375
$ws.Range("N110").Value = -8478.25
$ws.Range("H132").Value = 1880.9474
$ws.Range("I132").Value = 1807.7059
$ws.Range("J132").Value = 2503.5
$ws.Range("K132").Value = 5423.1177
$ws.Range("L132").Value = 7510.5
$ws.Range("M132").Value = -2893.1177
$ws.Range("N132").Value = -12570.5
$ws.Range("H136").Value = 3546.0625
$ws.Range("I136").Value = 3619.0908
$ws.Range("J136").Value = 3385.4
$ws.Range("K136").Value = 10857.2724
$ws.Range("L136").Value = 10156.2
$ws.Range("M136").Value = -8307.2724
$ws.Range("N136").Value = -15256.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5808.35
$ws.Range("J134").Value = 3468.8333
$ws.Range("L134").Value = 10406.4999
$ws.Range("N134").Value = -15476.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3333883.5
$ws.Range("I22").Value = 884
$ws.Range("K22").Value = 884
$ws.Range("M22").Value = -534
$ws.Range("H58").Value = 2849.4167
$ws.Range("I58").Value = 1665.8334
$ws.Range("K58").Value = 1665.8334
$ws.Range("M58").Value = -1462.8334
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H134").Value = 2189.0454
$ws.Range("I134").Value = 2342.2222
$ws.Range("K134").Value = 7026.6666
$ws.Range("M134").Value = -4491.6666
$ws.Range("H136").Value = 2849.4167
$ws.Range("I136").Value = 1665.8334
$ws.Range("K136").Value = 4997.5002
$ws.Range("M136").Value = -2447.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 19736.777
$ws.Range("I7").Value = 23400.285
$ws.Range("J7").Value = 6914.5
$ws.Range("K7").Value = 70200.855
$ws.Range("L7").Value = 20743.5
$ws.Range("M7").Value = -70088.855
$ws.Range("N7").Value = -20967.5
$ws.Range("H12").Value = 488.1111
$ws.Range("J12").Value = 399.25
$ws.Range("L12").Value = 1197.75
$ws.Range("N12").Value = -1543.75
$ws.Range("H131").Value = 1872.875
$ws.Range("I131").Value = 1749.5
$ws.Range("K131").Value = 5248.5
$ws.Range("M131").Value = -208.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 3000
$ws.Range("J35").Value = 3000
$ws.Range("L35").Value = 3000
$ws.Range("N35").Value = -3596
$ws.Range("H46").Value = 22437.223
$ws.Range("I46").Value = 41
$ws.Range("J46").Value = 25236.75
$ws.Range("K46").Value = 41
$ws.Range("L46").Value = 25236.75
$ws.Range("M46").Value = 115
$ws.Range("N46").Value = -25548.75
$ws.Range("H102").Value = 2299.6667
$ws.Range("I102").Value = 2128.2856
$ws.Range("J102").Value = 2899.5
$ws.Range("K102").Value = 2128.2856
$ws.Range("L102").Value = 2899.5
$ws.Range("M102").Value = -506.2856000000002
$ws.Range("N102").Value = -6143.5
$ws.Range("H113").Value = 2499.75
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 2999.6667
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 2999.6667
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -7339.6667
$ws.Range("H122").Value = 2331.5334
$ws.Range("I122").Value = 2014.6666
$ws.Range("K122").Value = 6043.9998
$ws.Range("M122").Value = -3593.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1799.8889
$ws.Range("I46").Value = 1108
$ws.Range("J46").Value = 3183.6667
$ws.Range("K46").Value = 1108
$ws.Range("L46").Value = 3183.6667
$ws.Range("M46").Value = -920
$ws.Range("N46").Value = -3559.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 39900
$ws.Range("I70").Value = 39900
$ws.Range("K70").Value = 39900
$ws.Range("M70").Value = -39585
$ws.Range("H73").Value = 39900
$ws.Range("I73").Value = 39900
$ws.Range("K73").Value = 39900
$ws.Range("M73").Value = -38808
$ws.Range("H93").Value = 43919.5
$ws.Range("J93").Value = 43919.5
$ws.Range("L93").Value = 43919.5
$ws.Range("N93").Value = -48911.5

